$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.970.44'
$ws.Range("E2").Value = '  +2.85%  '

$ws.Range("D3").Value = '1.861.22'
$ws.Range("E3").Value = '  +2.03%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.08'
$ws.Range("E5").Value = '  +1.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6388'
$ws.Range("E6").Value = '  +3.67%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2996'
$ws.Range("E8").Value = '  +3.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07483'
$ws.Range("E9").Value = '  +1.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.54'
$ws.Range("E10").Value = '  +6.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07685'
$ws.Range("E11").Value = '  +0.27%  '

$ws.Range("D12").Value = '1.871.30'
$ws.Range("E12").Value = '  +2.38%  '

$ws.Range("E13").Value = '  +2.17%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6908'
$ws.Range("E14").Value = '  +4.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '84.23'
$ws.Range("E15").Value = '  +2.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009458'
$ws.Range("E16").Value = '  +5.85%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.078'
$ws.Range("E17").Value = '  +4.22%  '

$ws.Range("D18").Value = '29.946.19'
$ws.Range("E18").Value = '  +2.93%  '

$ws.Range("D19").Value = '2.123.55'
$ws.Range("E19").Value = '  +2.58%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '239.79'
$ws.Range("E20").Value = '  +1.15%  '

$ws.Range("E21").Value = '  +2.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9999'
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.416'
$ws.Range("E23").Value = '  +4.07%  '

$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.75'
$ws.Range("E25").Value = '  +1.26%  '

$ws.Range("E26").Value = '  +0.91%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.586'
$ws.Range("E27").Value = '  +1.87%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.99'
$ws.Range("E28").Value = '  +2.18%  '

$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06118'
$ws.Range("E29").Value = '  +10.24%  '

$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.506'
$ws.Range("E30").Value = '  +1.44%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.282'
$ws.Range("E31").Value = '  +6.67%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.137'
$ws.Range("E32").Value = '  +1.20%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.147'
$ws.Range("E33").Value = '  +1.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.892'
$ws.Range("E34").Value = '  +3.66%  '

$ws.Range("E35").Value = '  +3.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7327'
$ws.Range("E36").Value = '  -0.17%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.604'
$ws.Range("E37").Value = '  -0.14%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.866'
$ws.Range("E38").Value = '  +1.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01805'
$ws.Range("E39").Value = '  +2.74%  '

$ws.Range("D40").Value = '1.224.74'
$ws.Range("E40").Value = '  +1.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9314'
$ws.Range("E41").Value = '  +3.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.276'
$ws.Range("E42").Value = '  -0.92%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.17%  '

$ws.Range("E44").Value = '  +2.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '102.20'
$ws.Range("E45").Value = '  +0.91%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '66.38'
$ws.Range("E46").Value = '  +2.88%  '

$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5103'
$ws.Range("E47").Value = '  +0.48%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000120'
$ws.Range("E48").Value = '  -2.20%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.349'
$ws.Range("E49").Value = '  +3.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4103'
$ws.Range("E50").Value = '  +2.51%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1144'
$ws.Range("E51").Value = '  +3.28%  '
